# Check.xlsx edit — "add tests and methods"
# Adds a new test-result row for the "Добавления комментария к заявке" test
# case: marks it automated (column D) and records the observed result in
# column E ("Время комментария не совпадает с текущим на устройстве").
# Also widens column E to fit the new longer text and updates the active
# selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 25 ("Добавления комментария к заявке") as automated — copy the
# format from the D18 cell (same "V" style used throughout column D) and set
# the value.
$ws.Range("D18").Copy()
[void]$ws.Range("D25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D25").Value2 = "V"

# Record the test result for that row.
$ws.Range("E25").Value2 = "Время комментария не совпадает с текущим на устройстве"

# Widen column E so the new, longer result text fits.
$ws.Columns.Item(5).ColumnWidth = 54.08

# Update the active selection to match the author's last cursor position.
[void]$ws.Range("E23").Select()
